$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.165.61'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '1.839.75'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.93'
$ws.Range("E5").Value = '  -0.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6269'
$ws.Range("E6").Value = '  -1.52%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07522'
$ws.Range("E8").Value = '  -0.84%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2944'
$ws.Range("E9").Value = '  -0.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.34'
$ws.Range("E10").Value = '  +1.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07710'

$ws.Range("D12").Value = '1.844.65'
$ws.Range("E12").Value = '  +0.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.023'
$ws.Range("E13").Value = '  +0.20%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6779'
$ws.Range("E14").Value = '  +0.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.22'
$ws.Range("E15").Value = '  -0.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009279'
$ws.Range("E16").Value = '  -4.25%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.975'
$ws.Range("E17").Value = '  -2.61%  '

$ws.Range("D18").Value = '29.146.64'
$ws.Range("E18").Value = '  +0.08%  '

$ws.Range("D19").Value = '2.090.04'
$ws.Range("E19").Value = '  -0.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '233.03'
$ws.Range("E20").Value = '  +2.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.71'
$ws.Range("E21").Value = '  +0.79%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.173'
$ws.Range("E23").Value = '  -0.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.38'
$ws.Range("E25").Value = '  -0.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1404'
$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.553'
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.96'
$ws.Range("E28").Value = '  -0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.496'
$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.192'
$ws.Range("E30").Value = '  +1.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.159'
$ws.Range("E31").Value = '  +1.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05582'
$ws.Range("E32").Value = '  +3.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.209'
$ws.Range("E33").Value = '  +0.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7550'
$ws.Range("E34").Value = '  +0.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.857'
$ws.Range("E35").Value = '  -0.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.150'
$ws.Range("E36").Value = '  +0.46%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.661'
$ws.Range("E37").Value = '  -0.16%  '

$ws.Range("D38").Value = '1.242.75'
$ws.Range("E38").Value = '  +0.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.770'
$ws.Range("E39").Value = '  +0.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01788'
$ws.Range("E40").Value = '  -0.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.614'
$ws.Range("E41").Value = '  -0.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9000'
$ws.Range("E42").Value = '  -0.85%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.57'
$ws.Range("E44").Value = '  +0.19%  '

$ws.Range("D45").Value = '1.991.30'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.79'
$ws.Range("E46").Value = '  +2.50%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5090'
$ws.Range("E47").Value = '  -0.63%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000120'
$ws.Range("E48").Value = '  -2.48%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4097'
$ws.Range("E49").Value = '  +0.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.088'
$ws.Range("E50").Value = '  +0.08%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07266'
$ws.Range("E51").Value = '  +8.81%  '
